# Adds a new "Pause" textbox to slide 1, matching the existing
# "Here" textbox's formatting (yellow fill, 40pt autosize text box).
#
# Shape.Left/Top/Width/Height are backed by single-precision (Single)
# floats in the PowerPoint object model, which truncates the EMU value
# we actually want. EmuToPtsExact searches for a point value that,
# once round-tripped through that Single, reproduces the exact target
# EMU on save.
function EmuToPtsExact([double]$targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -lt 20000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $f = [float]$candidate
        $emu = [math]::Floor([double]$f * 12700.0)
        if ($emu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the existing yellow "Here" textbox (id 8) to copy its formatting.
$here = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 7") {
        $here = $sh
    }
}

# New shapes are assigned the lowest unused shape id on the slide (3
# and 7 are free here). Temporarily occupy those slots so the
# duplicated shape we keep lands on id 9, matching the target deck.
$filler1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$filler2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)

$dupRange = $here.Duplicate()
$pause = $dupRange.Item(1)

$filler1.Delete()
$filler2.Delete()

$pause.Name = "TextBox 8"
$pause.Left = EmuToPtsExact 6985597
$pause.Top = EmuToPtsExact 368715
$pause.Width = EmuToPtsExact 1420732
$pause.Height = EmuToPtsExact 707886
$pause.TextFrame.TextRange.Text = "Pause"
